$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "Pondicherry" text to "Pondicherry -???"
$ws.Range("B1").Value = "Pondicherry -???"

# Add a new column C: duplicate of B1's new site name, and D2 = 3 (new site id)
$ws.Range("C1").Value = "Pondicherry -???"
$ws.Range("D2").Value = 3

# Widen columns B and C to fit new content (16 characters, matching the
# author's manual column resize in the edited workbook)
$ws.Columns.Item(2).ColumnWidth = 15.1666666666667
$ws.Columns.Item(3).ColumnWidth = 15.1666666666667

# Move the active selection to D4 as in the edited workbook
$ws.Range("D4").Select() | Out-Null
